$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original row 2 values (before any changes), since row 2's data
# wraps around to become the new row 6.
$origD2 = $ws.Range("D2").Value2
$origK2 = $ws.Range("K2").Value2
$origL2 = $ws.Range("L2").Value2
$origM2 = $ws.Range("M2").Value2
$origN2 = $ws.Range("N2").Value2
$origO2 = $ws.Range("O2").Value2
$origP2 = $ws.Range("P2").Value2
$origQ2 = $ws.Range("Q2").Value2
$origR2 = $ws.Range("R2").Value2
$origS2 = $ws.Range("S2").Value2
$origT2 = $ws.Range("T2").Value2

# Shift rows 3..6 up into rows 2..5 (columns D, K, L, M, N, O, P, Q, R, S, T)
for ($r = 2; $r -le 5; $r++) {
    $src = $r + 1
    $ws.Range("D$r").Value2 = $ws.Range("D$src").Value2
    $ws.Range("K$r").Value2 = $ws.Range("K$src").Value2
    $ws.Range("L$r").Value2 = $ws.Range("L$src").Value2
    $ws.Range("M$r").Value2 = $ws.Range("M$src").Value2
    $ws.Range("N$r").Value2 = $ws.Range("N$src").Value2
    $ws.Range("O$r").Value2 = $ws.Range("O$src").Value2
    $ws.Range("P$r").Value2 = $ws.Range("P$src").Value2
    $ws.Range("Q$r").Value2 = $ws.Range("Q$src").Value2
    $ws.Range("R$r").Value2 = $ws.Range("R$src").Value2
    $ws.Range("S$r").Value2 = $ws.Range("S$src").Value2
    $ws.Range("T$r").Value2 = $ws.Range("T$src").Value2
}

# Wrap original row 2 values into row 6
$ws.Range("D6").Value2 = $origD2
$ws.Range("K6").Value2 = $origK2
$ws.Range("L6").Value2 = $origL2
$ws.Range("M6").Value2 = $origM2
$ws.Range("N6").Value2 = $origN2
$ws.Range("O6").Value2 = $origO2
$ws.Range("P6").Value2 = $origP2
$ws.Range("Q6").Value2 = $origQ2
$ws.Range("R6").Value2 = $origR2
$ws.Range("S6").Value2 = $origS2
$ws.Range("T6").Value2 = $origT2
